$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.499.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4688"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07743"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.885.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7286"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.173"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.493.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007472"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.126.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.257"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.257"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.080"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.892"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09700"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.468"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.277"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.125"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04855"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6926"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.720"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01887"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.827"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.181"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.008"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4246"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8233"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.487"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.969"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "916.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  +1.85%  "
